$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 68140.28999999999
$ws.Range("J75").Value = 68140.28999999999
$ws.Range("L75").Value = 68140.28999999999
$ws.Range("N75").Value = -70012.28999999999
$ws.Range("H78").Value = 68140.28999999999
$ws.Range("J78").Value = 68140.28999999999
$ws.Range("L78").Value = 204420.87
$ws.Range("N78").Value = -213780.87
$ws.Range("H101").Value = 757.6429000000001
$ws.Range("J101").Value = 1287.7142
$ws.Range("L101").Value = 3863.1426
$ws.Range("N101").Value = -7107.142599999999
$ws.Range("H106").Value = 8371.682000000001
$ws.Range("I106").Value = 8371.682000000001
$ws.Range("K106").Value = 8371.682000000001
$ws.Range("M106").Value = -7740.682000000001
$ws.Range("H107").Value = 1428.375
$ws.Range("I107").Value = 920.6667
$ws.Range("K107").Value = 920.6667
$ws.Range("M107").Value = 999.3333
$ws.Range("H113").Value = 3144.3333
$ws.Range("I113").Value = 3081.4375
$ws.Range("J113").Value = 3270.125
$ws.Range("K113").Value = 3081.4375
$ws.Range("L113").Value = 3270.125
$ws.Range("M113").Value = 172.5625
$ws.Range("N113").Value = -9778.125
$ws.Range("H121").Value = 2235.8333
$ws.Range("J121").Value = 2235.8333
$ws.Range("L121").Value = 6707.499899999999
$ws.Range("N121").Value = -10201.4999
$ws.Range("H132").Value = 2001.22
$ws.Range("I132").Value = 2090.0444
$ws.Range("J132").Value = 1201.8
$ws.Range("K132").Value = 6270.1332
$ws.Range("L132").Value = 3605.4
$ws.Range("M132").Value = -3740.1332
$ws.Range("N132").Value = -8665.4
$ws.Range("H137").Value = 23813440
$ws.Range("I137").Value = 38464788
$ws.Range("J137").Value = 4998.75
$ws.Range("K137").Value = 115394364
$ws.Range("L137").Value = 14996.25
$ws.Range("M137").Value = -115391814
$ws.Range("N137").Value = -20096.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 1000000000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 1000000000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 1000000000
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -1000000616
$ws.Range("H32").Value = 8337.691999999999
$ws.Range("I32").Value = 7968.3784
$ws.Range("K32").Value = 7968.3784
$ws.Range("M32").Value = -7681.3784
$ws.Range("H74").Value = 1077.4762
$ws.Range("I74").Value = 923.7222
$ws.Range("K74").Value = 923.7222
$ws.Range("M74").Value = -49.72220000000004
$ws.Range("H77").Value = 1077.4762
$ws.Range("I77").Value = 923.7222
$ws.Range("K77").Value = 4618.611
$ws.Range("M77").Value = -250.6109999999999
$ws.Range("H88").Value = 2774.6
$ws.Range("I88").Value = 2431.5715
$ws.Range("J88").Value = 2959.3076
$ws.Range("K88").Value = 2431.5715
$ws.Range("L88").Value = 2959.3076
$ws.Range("M88").Value = -2025.5715
$ws.Range("N88").Value = -3771.3076
$ws.Range("H91").Value = 2774.6
$ws.Range("I91").Value = 2431.5715
$ws.Range("J91").Value = 2959.3076
$ws.Range("K91").Value = 2431.5715
$ws.Range("L91").Value = 2959.3076
$ws.Range("M91").Value = -1027.5715
$ws.Range("N91").Value = -5767.3076
$ws.Range("H122").Value = 2094.9607
$ws.Range("I122").Value = 2023.2046
$ws.Range("J122").Value = 2546
$ws.Range("K122").Value = 6069.6138
$ws.Range("L122").Value = 7638
$ws.Range("M122").Value = -3619.6138
$ws.Range("N122").Value = -12538
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1341.8636
$ws.Range("I86").Value = 1012.4375
$ws.Range("J86").Value = 2220.3333
$ws.Range("K86").Value = 1012.4375
$ws.Range("L86").Value = 2220.3333
$ws.Range("M86").Value = 110.5625
$ws.Range("N86").Value = -4466.3333
$ws.Range("H89").Value = 1341.8636
$ws.Range("I89").Value = 1012.4375
$ws.Range("J89").Value = 2220.3333
$ws.Range("K89").Value = 5062.1875
$ws.Range("L89").Value = 11101.6665
$ws.Range("M89").Value = 553.8125
$ws.Range("N89").Value = -22333.6665
$ws.Range("H107").Value = 3094.2222
$ws.Range("I107").Value = 3507.85
$ws.Range("J107").Value = 1912.4286
$ws.Range("K107").Value = 3507.85
$ws.Range("L107").Value = 1912.4286
$ws.Range("M107").Value = -1587.85
$ws.Range("N107").Value = -5752.4286
$ws.Range("H140").Value = 190999.4
$ws.Range("J140").Value = 190999.4
$ws.Range("L140").Value = 190999.4
$ws.Range("N140").Value = -201359.4
$ws.Range("H141").Value = 191654.9
$ws.Range("I141").Value = 186475
$ws.Range("J141").Value = 192172.9
$ws.Range("K141").Value = 186475
$ws.Range("L141").Value = 192172.9
$ws.Range("M141").Value = -181295
$ws.Range("N141").Value = -202532.9
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 260.25
$ws.Range("I7").Value = 246.33333
$ws.Range("J7").Value = 268.6
$ws.Range("K7").Value = 246.33333
$ws.Range("L7").Value = 268.6
$ws.Range("M7").Value = -133.33333
$ws.Range("N7").Value = -494.6
$ws.Range("H31").Value = 24392950
$ws.Range("I31").Value = 27780014
$ws.Range("J31").Value = 6088.8
$ws.Range("K31").Value = 27780014
$ws.Range("L31").Value = 6088.8
$ws.Range("M31").Value = -27779719
$ws.Range("N31").Value = -6678.8
$ws.Range("H34").Value = 24392950
$ws.Range("I34").Value = 27780014
$ws.Range("J34").Value = 6088.8
$ws.Range("K34").Value = 27780014
$ws.Range("L34").Value = 6088.8
$ws.Range("M34").Value = -27779812
$ws.Range("N34").Value = -6492.8
$ws.Range("H37").Value = 4500
$ws.Range("J37").Value = 4500
$ws.Range("L37").Value = 4500
$ws.Range("N37").Value = -4714
$ws.Range("H58").Value = 1836.174
$ws.Range("I58").Value = 1172.5883
$ws.Range("J58").Value = 3716.3333
$ws.Range("K58").Value = 1172.5883
$ws.Range("L58").Value = 3716.3333
$ws.Range("M58").Value = -969.5882999999999
$ws.Range("N58").Value = -4122.3333
$ws.Range("H99").Value = 14839.792
$ws.Range("I99").Value = 9468
$ws.Range("K99").Value = 9468
$ws.Range("M99").Value = -7970
$ws.Range("H103").Value = 33262.582
$ws.Range("I103").Value = 12862.667
$ws.Range("J103").Value = 53662.5
$ws.Range("K103").Value = 12862.667
$ws.Range("L103").Value = 53662.5
$ws.Range("M103").Value = -11690.667
$ws.Range("N103").Value = -56006.5
$ws.Range("H126").Value = 14839.792
$ws.Range("I126").Value = 9468
$ws.Range("K126").Value = 28404
$ws.Range("M126").Value = -25934
$ws.Range("H132").Value = 1243.1364
$ws.Range("I132").Value = 1087.7
$ws.Range("K132").Value = 3263.1
$ws.Range("M132").Value = -733.1000000000004
$ws.Range("H136").Value = 1836.174
$ws.Range("I136").Value = 1172.5883
$ws.Range("J136").Value = 3716.3333
$ws.Range("K136").Value = 3517.7649
$ws.Range("L136").Value = 11148.9999
$ws.Range("M136").Value = -967.7648999999997
$ws.Range("N136").Value = -16248.9999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 27787776
$ws.Range("I80").Value = 33344532
$ws.Range("K80").Value = 100033596
$ws.Range("M80").Value = -100032660
$ws.Range("H83").Value = 27787776
$ws.Range("I83").Value = 33344532
$ws.Range("K83").Value = 300100788
$ws.Range("M83").Value = -300096108
$ws.Range("H95").Value = 33333
$ws.Range("J95").Value = 33333
$ws.Range("L95").Value = 99999
$ws.Range("N95").Value = -104117
$ws.Range("I107").Value = 3724.5
$ws.Range("J107").Value = 4136810.2
$ws.Range("K107").Value = 11173.5
$ws.Range("L107").Value = 12410430.6
$ws.Range("M107").Value = -9253.5
$ws.Range("N107").Value = -12414270.6
$ws.Range("H132").Value = 3541.2727
$ws.Range("J132").Value = 3741.2727
$ws.Range("L132").Value = 33671.4543
$ws.Range("N132").Value = -38731.4543
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 78999.91
$ws.Range("J104").Value = 77249.75
$ws.Range("L104").Value = 77249.75
$ws.Range("N104").Value = -84237.75
$ws.Range("H122").Value = 947416.4
$ws.Range("I122").Value = 1103849.1
$ws.Range("J122").Value = 8820
$ws.Range("K122").Value = 3311547.3
$ws.Range("L122").Value = 26460
$ws.Range("M122").Value = -3309097.3
$ws.Range("N122").Value = -31360
$ws.Range("H132").Value = 2939755.5
$ws.Range("I132").Value = 2668.2036
$ws.Range("K132").Value = 8004.610799999999
$ws.Range("M132").Value = -5474.610799999999
$ws.Range("H136").Value = 10234.429
$ws.Range("J136").Value = 10234.429
$ws.Range("L136").Value = 30703.287
$ws.Range("N136").Value = -35803.287
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 66004000
$ws.Range("I22").Value = 132000000
$ws.Range("J22").Value = 8000
$ws.Range("K22").Value = 132000000
$ws.Range("L22").Value = 8000
$ws.Range("M22").Value = -131999705
$ws.Range("N22").Value = -8590
$ws.Range("H27").Value = 66004000
$ws.Range("I27").Value = 132000000
$ws.Range("J27").Value = 8000
$ws.Range("K27").Value = 132000000
$ws.Range("L27").Value = 8000
$ws.Range("M27").Value = -131999893
$ws.Range("N27").Value = -8214
$ws.Range("H46").Value = 1243.9445
$ws.Range("I46").Value = 1062.8
$ws.Range("J46").Value = 1313.6154
$ws.Range("K46").Value = 1062.8
$ws.Range("L46").Value = 1313.6154
$ws.Range("M46").Value = -874.8
$ws.Range("N46").Value = -1689.6154
$ws.Range("H57").Value = 29386.666
$ws.Range("I57").Value = 29386.666
$ws.Range("K57").Value = 29386.666
$ws.Range("M57").Value = -28820.666
$ws.Range("H136").Value = 4073.0386
$ws.Range("I136").Value = 3359.1365
$ws.Range("K136").Value = 10077.4095
$ws.Range("M136").Value = -7527.4095
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 107508.664
$ws.Range("J56").Value = 107508.664
$ws.Range("L56").Value = 107508.664
$ws.Range("N56").Value = -108936.664
$ws.Range("H105").Value = 44998
$ws.Range("J105").Value = 44998
$ws.Range("L105").Value = 44998
$ws.Range("N105").Value = -51986
$ws.Range("H136").Value = 207312.95
$ws.Range("I136").Value = 3662.8975
$ws.Range("J136").Value = 1001548.2
$ws.Range("K136").Value = 10988.6925
$ws.Range("L136").Value = 3004644.6
$ws.Range("M136").Value = -8438.692500000001
$ws.Range("N136").Value = -3009744.6
$ws.Range("H139").Value = 220715
$ws.Range("J139").Value = 220715
$ws.Range("L139").Value = 220715
$ws.Range("N139").Value = -230995
